$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 144.77777
$ws.Range("I4").Value = 92.166664
$ws.Range("K4").Value = 92.166664
$ws.Range("M4").Value = 21.833336
$ws.Range("H32").Value = 10349.6
$ws.Range("J32").Value = 10583.333
$ws.Range("L32").Value = 10583.333
$ws.Range("N32").Value = -11235.333
$ws.Range("H41").Value = 1426.125
$ws.Range("I41").Value = 1137.5
$ws.Range("J41").Value = 1714.75
$ws.Range("K41").Value = 1137.5
$ws.Range("L41").Value = 1714.75
$ws.Range("M41").Value = -697.5
$ws.Range("N41").Value = -2594.75
$ws.Range("H64").Value = 6496.6665
$ws.Range("I64").Value = 6593.6665
$ws.Range("J64").Value = 6399.6665
$ws.Range("K64").Value = 6593.6665
$ws.Range("L64").Value = 6399.6665
$ws.Range("M64").Value = -6345.6665
$ws.Range("N64").Value = -6895.6665
$ws.Range("H67").Value = 6496.6665
$ws.Range("I67").Value = 6593.6665
$ws.Range("J67").Value = 6399.6665
$ws.Range("K67").Value = 6593.6665
$ws.Range("L67").Value = 6399.6665
$ws.Range("M67").Value = -5735.6665
$ws.Range("N67").Value = -8115.6665
$ws.Range("H74").Value = 9220.5
$ws.Range("I74").Value = 9191
$ws.Range("K74").Value = 9191
$ws.Range("M74").Value = -8255
$ws.Range("H76").Value = 8331.75
$ws.Range("I76").Value = 8566.286
$ws.Range("J76").Value = 8003.4
$ws.Range("K76").Value = 8566.286
$ws.Range("L76").Value = 8003.4
$ws.Range("M76").Value = -8251.286
$ws.Range("N76").Value = -8633.4
$ws.Range("H77").Value = 9220.5
$ws.Range("I77").Value = 9191
$ws.Range("K77").Value = 45955
$ws.Range("M77").Value = -41275
$ws.Range("H79").Value = 8331.75
$ws.Range("I79").Value = 8566.286
$ws.Range("J79").Value = 8003.4
$ws.Range("K79").Value = 8566.286
$ws.Range("L79").Value = 8003.4
$ws.Range("M79").Value = -7474.286
$ws.Range("N79").Value = -10187.4
$ws.Range("H106").Value = 4031.25
$ws.Range("I106").Value = 4125
$ws.Range("J106").Value = 3750
$ws.Range("K106").Value = 4125
$ws.Range("L106").Value = 3750
$ws.Range("M106").Value = -3494
$ws.Range("N106").Value = -5012
$ws.Range("H137").Value = 12672.452
$ws.Range("I137").Value = 3254.7827
$ws.Range("K137").Value = 9764.348100000001
$ws.Range("M137").Value = -7214.348100000001
$ws.Range("H141").Value = 5051.4287
$ws.Range("J141").Value = 5037.5
$ws.Range("L141").Value = 15112.5
$ws.Range("N141").Value = -25472.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 214.83333
$ws.Range("J4").Value = 250
$ws.Range("L4").Value = 250
$ws.Range("N4").Value = -482
$ws.Range("H5").Value = 168
$ws.Range("J5").Value = 202
$ws.Range("L5").Value = 202
$ws.Range("N5").Value = -426
$ws.Range("H32").Value = 4249.5
$ws.Range("I32").Value = 1649
$ws.Range("K32").Value = 1649
$ws.Range("M32").Value = -1362
$ws.Range("H45").Value = 1512.5264
$ws.Range("I45").Value = 1492.6
$ws.Range("K45").Value = 1492.6
$ws.Range("M45").Value = -1115.6
$ws.Range("H97").Value = 2049.05
$ws.Range("I97").Value = 2077.9473
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 2077.9473
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -1581.9473
$ws.Range("N97").Value = -2492
$ws.Range("H132").Value = 456947.4
$ws.Range("I132").Value = 2472.362
$ws.Range("J132").Value = 2484605.2
$ws.Range("K132").Value = 7417.086
$ws.Range("L132").Value = 7453815.600000001
$ws.Range("M132").Value = -4887.086
$ws.Range("N132").Value = -7458875.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 168
$ws.Range("J4").Value = 202
$ws.Range("L4").Value = 202
$ws.Range("N4").Value = -432
$ws.Range("H94").Value = 479.07144
$ws.Range("I94").Value = 449.16
$ws.Range("J94").Value = 728.3333
$ws.Range("K94").Value = 449.16
$ws.Range("L94").Value = 728.3333
$ws.Range("M94").Value = 1.839999999999975
$ws.Range("N94").Value = -1630.3333
$ws.Range("H99").Value = 2574.3333
$ws.Range("I99").Value = 2300.238
$ws.Range("J99").Value = 4493
$ws.Range("K99").Value = 2300.238
$ws.Range("L99").Value = 4493
$ws.Range("M99").Value = -802.2379999999998
$ws.Range("N99").Value = -7489

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H94").Value = 13350.429
$ws.Range("J94").Value = 13628.556
$ws.Range("L94").Value = 13628.556
$ws.Range("N94").Value = -14530.556

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 176247.5
$ws.Range("J37").Value = 176247.5
$ws.Range("L37").Value = 528742.5
$ws.Range("N37").Value = -528966.5
$ws.Range("H122").Value = 15372287
$ws.Range("J122").Value = 4048601.5
$ws.Range("L122").Value = 36437413.5
$ws.Range("N122").Value = -36442313.5
$ws.Range("H131").Value = 2432.8586
$ws.Range("I131").Value = 10748.667
$ws.Range("J131").Value = 1601.2778
$ws.Range("K131").Value = 32246.001
$ws.Range("L131").Value = 4803.8334
$ws.Range("M131").Value = -27206.001
$ws.Range("N131").Value = -14883.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 374.9
$ws.Range("J22").Value = 447.25
$ws.Range("L22").Value = 447.25
$ws.Range("N22").Value = -1037.25
$ws.Range("H27").Value = 374.9
$ws.Range("J27").Value = 447.25
$ws.Range("L27").Value = 447.25
$ws.Range("N27").Value = -661.25
$ws.Range("H82").Value = 2944
$ws.Range("I82").Value = 1375
$ws.Range("J82").Value = 3990
$ws.Range("K82").Value = 1375
$ws.Range("L82").Value = 3990
$ws.Range("M82").Value = -1014
$ws.Range("N82").Value = -4712
$ws.Range("H85").Value = 2944
$ws.Range("I85").Value = 1375
$ws.Range("J85").Value = 3990
$ws.Range("K85").Value = 1375
$ws.Range("L85").Value = 3990
$ws.Range("M85").Value = -127
$ws.Range("N85").Value = -6486
$ws.Range("H93").Value = 14552.875
$ws.Range("I93").Value = 14552.875
$ws.Range("K93").Value = 14552.875
$ws.Range("M93").Value = -13304.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 580
$ws.Range("I100").Value = 501.5
$ws.Range("J100").Value = 674.2
$ws.Range("K100").Value = 1003
$ws.Range("L100").Value = 1348.4
$ws.Range("M100").Value = -462
$ws.Range("N100").Value = -2430.4
$ws.Range("H136").Value = 261106.08
$ws.Range("I136").Value = 1748.8422
$ws.Range("J136").Value = 485096.4
$ws.Range("K136").Value = 5246.5266
$ws.Range("L136").Value = 1455289.2
$ws.Range("M136").Value = -2696.5266
$ws.Range("N136").Value = -1460389.2
